# Propuesta Inicial del proyecto.docx - apply commit changes
#
# The commit's OOXML diff is dominated by Word's automatic spell-checker
# marking loan-words / proper nouns with <w:proofErr> (spellStart/spellEnd)
# -- a purely cosmetic, invisible (squiggly-underline) marker that does not
# change the rendered text and has no Word object-model property to author
# directly. The one real, visible, user-authored edit in the diff is a
# yellow highlight applied to the "Una síntesis..." abstract-instructions
# paragraph. We apply that here.

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Una síntesis de lo hecho hasta el momento", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $para = $rng.Paragraphs(1).Range
    $para.HighlightColorIndex = 7
    Write-Host "Highlighted paragraph: $($para.Text)"
} else {
    Write-Host "WARNING: target paragraph not found"
}
